# Workbook was edited to record explicit "0" time spent values for the
# first task of each participation (previously left blank), and the
# sheet view scroll/selection state was reset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose F column (time_spent_till_switching_to_this_task) was empty
# and must now contain an explicit 0 value.
$rows = @(2,9,16,24,27,35,40,48,52,57,63,68,74,77,81,84,89,92,95,101,105,111,116,121,127,133,138,143)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Reset the view: scroll back to the top (no topLeftCell override) and
# move the active selection to G1.
$ws.Range("G1").Select()
